$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Rename the existing sheet and add the three new sheets needed,
# in the correct order, so sheetId/order match the target workbook.
# -----------------------------------------------------------------
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sales vs PO"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

# -----------------------------------------------------------------
# Sheet 1: "Sales vs PO"
# Columns: A=ds, B=y, C=Order Week, D=PO_Requested_Qty
# -----------------------------------------------------------------
$ws1.Cells.Item(1,1).Value = "ds"
$ws1.Cells.Item(1,2).Value = "y"
$ws1.Cells.Item(1,3).Value = "Order Week"
$ws1.Cells.Item(1,4).Value = "PO_Requested_Qty"

$sheet1Rows = @(
    @(45571, 0, 45565, 0),
    @(45578, 0, 45572, 0),
    @(45585, 0, 45579, 0),
    @(45592, 0, 45586, 0),
    @(45599, 8, 45593, 0),
    @(45606, 6, 45600, 0),
    @(45613, 3, 45607, 0),
    @(45620, 7, 45614, 0),
    @(45627, 7, 45621, 0),
    @(45634, 13, 45628, 0),
    @(45641, 10, 45635, 0),
    @(45648, 4, 45642, 0),
    @(45655, 6, 45649, 0)
)

$r = 2
foreach ($row in $sheet1Rows) {
    $aCell = $ws1.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws1.Cells.Item($r, 2).Value = $row[1]

    $cCell = $ws1.Cells.Item($r, 3)
    $cCell.Value = $row[2]
    $cCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws1.Cells.Item($r, 4).Value = $row[3]

    $r = $r + 1
}

# Apply the existing header style (bold, centered, top-aligned, thin
# border) to the new D1 header cell by copying formats from an
# already-styled header cell (A1) - this reuses the same style
# record instead of creating a new, slightly different one.
$ws1.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)

# -----------------------------------------------------------------
# Sheet 2: "Weekly Growth"
# Columns: A=ds, B=PO_Requested_Qty, C=Growth%
# -----------------------------------------------------------------
$ws2.Cells.Item(1,1).Value = "ds"
$ws2.Cells.Item(1,2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1,3).Value = "Growth%"

$sheet2Rows = @(
    @(45572, 240, 0),
    @(45586, 16, -93.33333333333333),
    @(45607, 16, 0)
)

$r = 2
foreach ($row in $sheet2Rows) {
    $aCell = $ws2.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]

    $r = $r + 1
}

$ws1.Range("A1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

# -----------------------------------------------------------------
# Sheet 3: "Volume Insights"
# Columns: A=Total_PO_Quantity, B=Average_PO_Quantity,
#          C=Max_PO_Quantity, D=Min_PO_Quantity
# -----------------------------------------------------------------
$ws3.Cells.Item(1,1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1,2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1,3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1,4).Value = "Min_PO_Quantity"

$ws3.Cells.Item(2,1).Value = 272
$ws3.Cells.Item(2,2).Value = 90.66666666666667
$ws3.Cells.Item(2,3).Value = 240
$ws3.Cells.Item(2,4).Value = 16

$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# -----------------------------------------------------------------
# Sheet 4: "Prediction Info"
# Columns: A=Predicted_Next_Week_PO_Quantity
# -----------------------------------------------------------------
$ws4.Cells.Item(1,1).Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Cells.Item(2,1).Value = 0

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

# Clear the clipboard "marching ants" marquee left from the copy ops.
$excel.CutCopyMode = $false

# Re-select A1 on sheet 1 and make it the active sheet, mirroring the
# original workbook's view state.
$ws1.Activate()
$ws1.Range("A1").Select()
